# Insert a new row at position 232 in Sheet1. This shifts the existing
# rows 232-350 down to 233-351 and leaves a blank row (with the D-column
# date style inherited from the row above) at 232 for the new weekly
# price entry (Pomelo / Start Ruby / Primera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new weekly record.
$ws.Range("A232").Value = 10
$ws.Range("B232").Value = "Vega Modelo de Temuco"
$ws.Range("C232").Value = "La Araucanía"
$ws.Range("D232").Value = 44917
$ws.Range("E232").Value = 9
$ws.Range("F232").Value = "Fruta"
$ws.Range("G232").Value = 100102
$ws.Range("H232").Value = "Cítricos"
$ws.Range("I232").Value = 100102006
$ws.Range("J232").Value = "Pomelo"
$ws.Range("K232").Value = "Start Ruby"
$ws.Range("L232").Value = "Primera"
$ws.Range("M232").Value = 100
$ws.Range("N232").Value = 14000
$ws.Range("O232").Value = 14000
$ws.Range("P232").Value = 14000
$ws.Range("Q232").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R232").Value = "Región de O'Higgins"
$ws.Range("S232").Value = 933
$ws.Range("T232").Value = 15
